# [MOSIP-14336] Updating Masterdata Utility
# Re-layout the individual_type master-data sheet: add a "lang_code"
# column, rename headers to snake_case, and add French-language rows
# for the existing Foreigner / Non-Foreigner individual types.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row --------------------------------------------------------
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "is_active"

# Give the new lang_code header cell the same direct formatting (bold
# font, thin border, centered alignment) already used by the other
# header cells, by copying the format from an existing header cell.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Data rows -----------------------------------------------------------
# Row 2: eng / FR / Foreigner / TRUE
$ws.Range("A2").Value = "eng"
$ws.Range("B2").Value = "FR"
$ws.Range("C2").Value = "Foreigner"
$ws.Range("D2").Value = $true

# Row 3: eng / NFR / Non-Foreigner / TRUE
$ws.Range("A3").Value = "eng"
$ws.Range("B3").Value = "NFR"
$ws.Range("C3").Value = "Non-Foreigner"
$ws.Range("D3").Value = $true

# Row 4: fra / FR / Étranger / TRUE
$ws.Range("A4").Value = "fra"
$ws.Range("B4").Value = "FR"
$ws.Range("C4").Value = "Étranger"
$ws.Range("D4").Value = $true

# Row 5 (new row): fra / NFR / Non-étranger / TRUE
$ws.Range("A5").Value = "fra"
$ws.Range("B5").Value = "NFR"
$ws.Range("C5").Value = "Non-étranger"
$ws.Range("D5").Value = $true

# The lang_code data cells (A2:A5) carry the same direct formatting as
# the header cells (s="1" in the target file) - copy it across.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A2:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0
